# New PO forecast model
# Updates the three worksheets (Weekly Quantity, Monthly Trend, PO Forecast)
# with refreshed forecast data: appends newly observed weekly/monthly rows
# and replaces the PO Forecast projection with the new model's output.

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append 4 new weekly rows (20-23)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

$weekly = @(
    @(20, 45662.99999999999, 68),
    @(21, 45669.99999999999, 20),
    @(22, 45676.99999999999, 24),
    @(23, 45683.99999999999, 44)
)

foreach ($row in $weekly) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws1.Cells.Item($r, 2).Value = $row[2]
}

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append 1 new monthly row (10)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Cells.Item(10, 1).Value = 45688.99999999999
$ws2.Cells.Item(10, 1).NumberFormat = $dateFormat
$ws2.Cells.Item(10, 2).Value = 156

# ---------------------------------------------------------------------
# Sheet 3: "PO Forecast" - new forecast model values for existing rows
# (2-19) and append new forecast rows (20-31)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PO Forecast")

$forecast = @(
    @(2, 45417.99999999999, 39),
    @(3, 45473.99999999999, 47),
    @(4, 45487.99999999999, 49),
    @(5, 45501.99999999999, 51),
    @(6, 45508.99999999999, 52),
    @(7, 45515.99999999999, 53),
    @(8, 45522.99999999999, 54),
    @(9, 45529.99999999999, 55),
    @(10, 45543.99999999999, 57),
    @(11, 45557.99999999999, 59),
    @(12, 45564.99999999999, 60),
    @(13, 45571.99999999999, 61),
    @(14, 45578.99999999999, 62),
    @(15, 45592.99999999999, 64),
    @(16, 45599.99999999999, 65),
    @(17, 45627.99999999999, 69),
    @(18, 45634.99999999999, 70),
    @(19, 45641.99999999999, 71),
    @(20, 45662.99999999999, 74),
    @(21, 45669.99999999999, 75),
    @(22, 45676.99999999999, 76),
    @(23, 45683.99999999999, 77),
    @(24, 45690.99999999999, 78),
    @(25, 45697.99999999999, 79),
    @(26, 45704.99999999999, 80),
    @(27, 45711.99999999999, 81),
    @(28, 45718.99999999999, 82),
    @(29, 45725.99999999999, 83),
    @(30, 45732.99999999999, 84),
    @(31, 45739.99999999999, 84)
)

foreach ($row in $forecast) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 1).NumberFormat = $dateFormat
    $ws3.Cells.Item($r, 2).Value = $row[2]
}
